# Applies updated crypto price/volume figures to columns D (Price) and E (Volume(1h))
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: cell address -> new text value.
# Values are written as literal text (matching the sheet's original inline-string
# cells), so numeric-looking prices are forced to Text format first and the style
# is reset back to Normal afterwards to avoid leaving a stray number format behind.
$updates = [ordered]@{
    'D2' = '26.002.30'
    'E2' = '  -0.27%  '
    'D3' = '1.745.81'
    'E3' = '  -0.19%  '
    'D4' = '1.001'
    'D5' = '248.92'
    'E5' = '  +6.25%  '
    'D6' = '1.001'
    'E6' = '  +0.24%  '
    'D7' = '0.5154'
    'E7' = '  -1.83%  '
    'D8' = '0.2761'
    'E8' = '  -0.94%  '
    'D9' = '0.06202'
    'E9' = '  +0.07%  '
    'D10' = '1.740.47'
    'E10' = '  -0.54%  '
    'D11' = '0.07244'
    'E11' = '  +0.96%  '
    'D12' = '15.20'
    'E12' = '  -1.19%  '
    'D13' = '0.6493'
    'E13' = '  +0.51%  '
    'D14' = '4.634'
    'D15' = '77.86'
    'E15' = '  -0.47%  '
    'D16' = '1.001'
    'E16' = '  +0.26%  '
    'D17' = '1.001'
    'E17' = '  +0.19%  '
    'D18' = '26.044.41'
    'E18' = '  +0.27%  '
    'E19' = '  +1.75%  '
    'E20' = '  +1.78%  '
    'D21' = '1.965.35'
    'E21' = '  -0.58%  '
    'D22' = '4.299'
    'E22' = '  -0.45%  '
    'D23' = '8.694'
    'E23' = '  -1.61%  '
    'D24' = '5.365'
    'E24' = '  +3.00%  '
    'D25' = '135.46'
    'E25' = '  -2.68%  '
    'D26' = '1.506'
    'E26' = '  -0.62%  '
    'E27' = '  -0.34%  '
    'D28' = '1.785'
    'E28' = '  -1.99%  '
    'D29' = '106.14'
    'E29' = '  +1.81%  '
    'D30' = '3.964'
    'E30' = '  +4.80%  '
    'D31' = '0.08268'
    'E31' = '  -0.94%  '
    'D32' = '3.678'
    'E32' = '  -0.46%  '
    'D33' = '0.04678'
    'D34' = '2.656'
    'E34' = '  +0.89%  '
    'D35' = '1.001'
    'E35' = '  -0.07%  '
    'D36' = '0.6231'
    'E36' = '  -0.96%  '
    'D37' = '2.740'
    'E37' = '  +1.13%  '
    'E38' = '  +0.27%  '
    'D39' = '1.935'
    'E39' = '  -0.46%  '
    'D40' = '1.0000'
    'E40' = '  +0.24%  '
    'D41' = '100.30'
    'E41' = '  +1.93%  '
    'D42' = '0.3881'
    'E42' = '  -0.67%  '
    'D43' = '0.7571'
    'E43' = '  +2.92%  '
    'D44' = '5.020'
    'E44' = '  -1.05%  '
    'D45' = '6.346'
    'E45' = '  +0.32%  '
    'E46' = '  -0.19%  '
    'D47' = '55.49'
    'E47' = '  +2.74%  '
    'D48' = '0.05232'
    'E48' = '  -2.41%  '
    'D49' = '30.71'
    'E49' = '  +0.98%  '
    'D50' = '7.628'
    'E50' = '  -1.06%  '
    'D51' = '0.3447'
    'E51' = '  -0.76%  '
}

foreach ($addr in $updates.Keys) {
    $value = $updates[$addr]
    $cell = $ws.Range($addr)
    if ($addr.StartsWith("D")) {
        # Column D holds price text; pre-format as Text so numeric-looking
        # strings (e.g. "1.001") are not auto-converted to a Double by COM.
        $cell.NumberFormat = "@"
        $cell.Value = $value
        $cell.Style = "Normal"
    } else {
        # Column E values (e.g. "  +0.24%  ") are never numeric, so a plain
        # assignment keeps them as text without touching the cell style.
        $cell.Value = $value
    }
}
